$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = -0.096431111944355052
$ws.Range("B1").Value = 0.096096837119489464
$ws.Range("A2").Value = -0.04998835119440237
$ws.Range("B2").Value = 0.048623992764555979
$ws.Range("A3").Value = 0.11260483263900412
$ws.Range("B3").Value = -0.11314374154746787
$ws.Range("A4").Value = -0.17084425943638948
$ws.Range("B4").Value = 0.16992783885789464
$ws.Range("A5").Value = -0.16392783951030232
$ws.Range("B5").Value = 0.16208588450520889
$ws.Range("A6").Value = -0.10938831666776538
$ws.Range("B6").Value = 0.10920833069037439
$ws.Range("A7").Value = -0.089208331497967919
$ws.Range("B7").Value = 0.088754974750147042
$ws.Range("A8").Value = -0.068754975567644649
$ws.Range("B8").Value = 0.06836257598900275
$ws.Range("A9").Value = -0.062362576690385474
$ws.Range("B9").Value = 0.062028468063493492
$ws.Range("A10").Value = -0.056028468776226248
$ws.Range("B10").Value = 0.05598338477007303
$ws.Range("A11").Value = -0.051483385470440624
$ws.Range("B11").Value = 0.051404338034540586
$ws.Range("A12").Value = -0.045404338751324325
$ws.Range("B12").Value = 0.045153933599204166
$ws.Range("A13").Value = -0.039153934328109763
$ws.Range("B13").Value = 0.039086011929055253
$ws.Range("A14").Value = -0.027086012716592833
$ws.Range("B14").Value = 0.027053538167981728
$ws.Range("A15").Value = -0.021053538903386126
$ws.Range("B15").Value = 0.021028013487283559
$ws.Range("A16").Value = -0.015028014225265229
$ws.Range("B16").Value = 0.015004216383139202
$ws.Range("A17").Value = -0.0090042171245769964
$ws.Range("B17").Value = 0.0089999992303013698
$ws.Range("A18").Value = -0.036108486605066048
$ws.Range("B18").Value = 0.036096262020315351
$ws.Range("A19").Value = -0.027096262691129969
$ws.Range("B19").Value = 0.027013245846056577
$ws.Range("A20").Value = -0.018013246522913917
$ws.Range("B20").Value = 0.018004248487798691
$ws.Range("A21").Value = -0.0090042491655539791
$ws.Range("B21").Value = 0.0089999993215839069
$ws.Range("A22").Value = -0.093936291586004472
$ws.Range("B22").Value = 0.093627412585815861
$ws.Range("A23").Value = -0.084627413267226004
$ws.Range("B23").Value = 0.084125277573972035
$ws.Range("A24").Value = -0.042125278568402535
$ws.Range("B24").Value = 0.041999999000113419
$ws.Range("A25").Value = -0.046697568513515364
$ws.Range("B25").Value = 0.046657750847323598
$ws.Range("A26").Value = -0.040657751524332042
$ws.Range("B26").Value = 0.040610958361362748
$ws.Range("A27").Value = -0.034610959039802047
$ws.Range("B27").Value = 0.034462324151417256
$ws.Range("A28").Value = -0.028462324835929032
$ws.Range("B28").Value = 0.02837620675455188
$ws.Range("A29").Value = -0.016376207497627249
$ws.Range("B29").Value = 0.016353586547738885
$ws.Range("A30").Value = 0.0036464126348625392
$ws.Range("B30").Value = -0.0037429277635916414
$ws.Range("A31").Value = 0.018190507889581653
$ws.Range("B31").Value = -0.018224957167815248
$ws.Range("A32").Value = 0.039224956342635231
$ws.Range("B32").Value = -0.039310346265988727
